$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Planilha_01")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Planilha_02"

$data = @(
    @("CAMPO_01","CAMPO_02","CAMPO_02"),
    @("VALOR_01_01","VALOR_02_1","VALOR_03_1"),
    @("VALOR_01_2","VALOR_02_2","VALOR_03_2"),
    @("VALOR_01_3","VALOR_02_3","VALOR_03_3"),
    @("VALOR_01_4","VALOR_02_4","VALOR_03_4"),
    @("VALOR_01_5","VALOR_02_5","VALOR_03_5"),
    @("VALOR_01_6","VALOR_02_6","VALOR_03_6"),
    @("VALOR_01_7","VALOR_02_7","VALOR_03_7"),
    @("VALOR_01_8","VALOR_02_8","VALOR_03_8"),
    @("VALOR_01_9","VALOR_02_9","VALOR_03_9"),
    @("VALOR_01_10","VALOR_02_10","VALOR_03_10"),
    @("VALOR_01_11","VALOR_02_11","VALOR_03_11"),
    @("VALOR_01_12","VALOR_02_12","VALOR_03_12"),
    @("VALOR_01_13","VALOR_02_13","VALOR_03_13")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
